$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Form Responses 1" survey rows 414-418 (5 new submissions).
# Style source cells (already present on the sheet) used purely to copy
# the correct cell style onto each new cell before writing its value:
#   A column -> date/time style   (template: A23, style id 3)
#   G column -> percentage style  (template: G23, style id 4)
#   everything else -> plain text/number style (template: B23, style id 2)
$dateStyleCell = "A23"
$pctStyleCell  = "G23"
$txtStyleCell  = "B23"

$newRows = @(
  @{ Row = 414; Cells = @(
      @{ Col = "A"; Value = 44246.629239525464; Style = $dateStyleCell },
      @{ Col = "B"; Value = "PK-Seutu (Helsinki, Espoo, Vantaa)"; Style = $txtStyleCell },
      @{ Col = "C"; Value = "36-40 v"; Style = $txtStyleCell },
      @{ Col = "D"; Value = "Nainen"; Style = $txtStyleCell },
      @{ Col = "E"; Value = 14; Style = $txtStyleCell },
      @{ Col = "F"; Value = "Työntekijä / palkollinen"; Style = $txtStyleCell },
      @{ Col = "G"; Value = 1; Style = $pctStyleCell },
      @{ Col = "H"; Value = "Senior consultant"; Style = $txtStyleCell },
      @{ Col = "I"; Value = "Pääosin tai kokonaan etätyö"; Style = $txtStyleCell },
      @{ Col = "J"; Value = 8500; Style = $txtStyleCell },
      @{ Col = "K"; Value = 100000; Style = $txtStyleCell },
      @{ Col = "L"; Value = "Kyllä"; Style = $txtStyleCell },
      @{ Col = "M"; Value = "Sulava"; Style = $txtStyleCell }
    ) },
  @{ Row = 415; Cells = @(
      @{ Col = "A"; Value = 44246.6346266551; Style = $dateStyleCell },
      @{ Col = "B"; Value = "Pori"; Style = $txtStyleCell },
      @{ Col = "C"; Value = "36-40 v"; Style = $txtStyleCell },
      @{ Col = "D"; Value = "Mies"; Style = $txtStyleCell },
      @{ Col = "E"; Value = 8; Style = $txtStyleCell },
      @{ Col = "F"; Value = "Työntekijä / palkollinen"; Style = $txtStyleCell },
      @{ Col = "G"; Value = 1; Style = $pctStyleCell },
      @{ Col = "H"; Value = "Tech Lead"; Style = $txtStyleCell },
      @{ Col = "I"; Value = "Pääosin tai kokonaan etätyö"; Style = $txtStyleCell },
      @{ Col = "J"; Value = 5080; Style = $txtStyleCell },
      @{ Col = "K"; Value = 65000; Style = $txtStyleCell },
      @{ Col = "L"; Value = "Ei"; Style = $txtStyleCell },
      @{ Col = "M"; Value = "Iso konsulttitalo"; Style = $txtStyleCell },
      @{ Col = "N"; Value = "Sijainti Pori, mutta etätöitä 100%. Varsinainen positio Tampere - Helsinki. Edut aika huonot, perusjutut. Työ itsessään aika masentavaa. Seuraavaksi varmaan freelance/yrittäjyys."; Style = $txtStyleCell }
    ) },
  @{ Row = 416; Cells = @(
      @{ Col = "A"; Value = 44246.641679224536; Style = $dateStyleCell },
      @{ Col = "B"; Value = "Tampere"; Style = $txtStyleCell },
      @{ Col = "C"; Value = "36-40 v"; Style = $txtStyleCell },
      @{ Col = "D"; Value = "Mies"; Style = $txtStyleCell },
      @{ Col = "E"; Value = 14; Style = $txtStyleCell },
      @{ Col = "F"; Value = "Työntekijä / palkollinen"; Style = $txtStyleCell },
      @{ Col = "G"; Value = 1; Style = $pctStyleCell },
      @{ Col = "H"; Value = "Ohjelmistotestaaja"; Style = $txtStyleCell },
      @{ Col = "I"; Value = "Pääosin tai kokonaan etätyö"; Style = $txtStyleCell },
      @{ Col = "J"; Value = 4100; Style = $txtStyleCell },
      @{ Col = "K"; Value = 55000; Style = $txtStyleCell },
      @{ Col = "L"; Value = "Kyllä"; Style = $txtStyleCell }
    ) },
  @{ Row = 417; Cells = @(
      @{ Col = "A"; Value = 44246.64923311342; Style = $dateStyleCell },
      @{ Col = "B"; Value = "Tampere"; Style = $txtStyleCell },
      @{ Col = "C"; Value = "26-30 v"; Style = $txtStyleCell },
      @{ Col = "D"; Value = "?"; Style = $txtStyleCell },
      @{ Col = "E"; Value = 7; Style = $txtStyleCell },
      @{ Col = "F"; Value = "Työntekijä / palkollinen"; Style = $txtStyleCell },
      @{ Col = "G"; Value = 1; Style = $pctStyleCell },
      @{ Col = "H"; Value = "Full-stack developer"; Style = $txtStyleCell },
      @{ Col = "I"; Value = "Noin 50/50 hybridimalli"; Style = $txtStyleCell },
      @{ Col = "J"; Value = 5550; Style = $txtStyleCell },
      @{ Col = "K"; Value = 69400; Style = $txtStyleCell },
      @{ Col = "L"; Value = "Kyllä"; Style = $txtStyleCell }
    ) },
  @{ Row = 418; Cells = @(
      @{ Col = "A"; Value = 44246.65296685185; Style = $dateStyleCell },
      @{ Col = "B"; Value = "PK-Seutu (Helsinki, Espoo, Vantaa)"; Style = $txtStyleCell },
      @{ Col = "C"; Value = "26-30 v"; Style = $txtStyleCell },
      @{ Col = "D"; Value = "Mies"; Style = $txtStyleCell },
      @{ Col = "E"; Value = 5; Style = $txtStyleCell },
      @{ Col = "F"; Value = "Työntekijä / palkollinen"; Style = $txtStyleCell },
      @{ Col = "G"; Value = 0.8; Style = $pctStyleCell },
      @{ Col = "H"; Value = "Full-stack/mobiili/design"; Style = $txtStyleCell },
      @{ Col = "I"; Value = "Pääosin tai kokonaan etätyö"; Style = $txtStyleCell },
      @{ Col = "J"; Value = 7000; Style = $txtStyleCell },
      @{ Col = "K"; Value = 90000; Style = $txtStyleCell },
      @{ Col = "L"; Value = "Kyllä"; Style = $txtStyleCell },
      @{ Col = "M"; Value = "Mavericks "; Style = $txtStyleCell }
    ) }
)

foreach ($r in $newRows) {
    foreach ($c in $r.Cells) {
        $target = "$($c.Col)$($r.Row)"
        $ws.Range($c.Style).Copy()
        $ws.Range($target).PasteSpecial(-4122)
        $ws.Range($target).Value = $c.Value
    }
}

